$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Z2").Value = "2025-11-02T02:03:40.604695"
$ws.Range("Z3").Value = "2025-11-02T02:03:40.605702"
$ws.Range("Z4").Value = "2025-11-02T02:03:40.605702"
$ws.Range("Z5").Value = "2025-11-02T02:03:40.605702"
$ws.Range("Z6").Value = "2025-11-02T02:03:40.605702"
$ws.Range("Z7").Value = "2025-11-02T02:03:40.605702"
$ws.Range("Z8").Value = "2025-11-02T02:03:40.606753"
$ws.Range("Z9").Value = "2025-11-02T02:03:40.606884"
$ws.Range("Z10").Value = "2025-11-02T02:03:40.606884"
$ws.Range("Z11").Value = "2025-11-02T02:03:40.606884"
$ws.Range("Z12").Value = "2025-11-02T02:03:40.609410"
$ws.Range("Z13").Value = "2025-11-02T02:03:40.609711"
$ws.Range("Z14").Value = "2025-11-02T02:03:40.609711"
$ws.Range("Z15").Value = "2025-11-02T02:03:40.610244"
$ws.Range("Z16").Value = "2025-11-02T02:03:40.610244"
$ws.Range("Z17").Value = "2025-11-02T02:03:40.610244"
$ws.Range("Z18").Value = "2025-11-02T02:03:40.610244"
$ws.Range("Z19").Value = "2025-11-02T02:03:40.610244"
$ws.Range("Z20").Value = "2025-11-02T02:03:40.610244"
$ws.Range("Z21").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z22").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z23").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z24").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z25").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z26").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z27").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z28").Value = "2025-11-02T02:03:40.611240"
$ws.Range("Z29").Value = "2025-11-02T02:03:40.611240"
